$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.338045954704285
$ws.Range("B1").Value = 2.371382474899292
$ws.Range("C1").Value = 2.907703638076782
$ws.Range("D1").Value = 3.348934412002563
$ws.Range("E1").Value = 1.720772624015808
